# Update on sprint backlog 09-04
# Fill in burndown values for the 2019-04-09 column (H) and carry them
# forward into the newly active column (I) for the active sprint tasks.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1. Iteration")

# Column I (04-09) gets the same "remaining work" values as column H for
# most rows (they were untouched previously).
$ws.Range("I11").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("I13").Value = 7
$ws.Range("I14").Value = 10
$ws.Range("I15").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("I19").Value = 0

# Row 20 previously had placeholder text values ("??") in H; replace with
# numeric 0 for both H and I.
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0

$ws.Range("I32").Value = 8
$ws.Range("I33").Value = 4

# Row 34 previously had placeholder text values ("???") in H; replace with
# numeric 0 for both H and I.
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0

# Update the view to reflect where the user was working when saving.
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("A11").Select()
